$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (Jaden Hardy) is replaced with a new player entry.
$ws.Range("A3").Value = "Ziaire Williams"
$ws.Range("B3").Value = "SG,SF"
$ws.Range("C3").Value = "Brooklyn Nets"
